# Fill in quality (E) and contact_angle (F) measurements for rows whose
# images have now been fully processed. The style used by existing
# quality/contact_angle cells (e.g. E2:F2) is copied onto the newly
# populated cells so the formatting matches the rest of the column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, quality, contact_angle
$data = @(
    ,@(11, 3.0, 61.13)
    ,@(12, 3.0, 63.67)
    ,@(201, 2.0, 58.62)
    ,@(202, 2.0, 61.56)
    ,@(203, 2.0, 61.47)
    ,@(204, 1.0, 57.11)
    ,@(205, 1.0, 59.25)
    ,@(206, 1.0, 57.59)
    ,@(207, 2.0, 53.59)
    ,@(208, 1.0, 64.75)
    ,@(209, 2.0, 53.39)
    ,@(210, 1.0, 64.19)
    ,@(211, 2.0, 63.97)
    ,@(212, 1.0, 60.65)
    ,@(213, 2.0, 58.68)
    ,@(214, 2.0, 58.61)
    ,@(215, 2.0, 53.31)
    ,@(216, 2.0, 65.28)
    ,@(217, 1.0, 60.51)
    ,@(218, 1.0, 62.15)
    ,@(219, 2.0, 57.66)
    ,@(220, 1.0, 70.9)
    ,@(221, 1.0, 62.73)
    ,@(222, 1.0, 62.39)
    ,@(223, 1.0, 69.67)
    ,@(224, 2.0, 59.34)
    ,@(225, 2.0, 68.0)
    ,@(226, 2.0, 68.06)
    ,@(227, 1.0, 76.21)
    ,@(228, 1.0, 58.58)
    ,@(229, 2.0, 54.56)
    ,@(230, 1.0, 31.34)
    ,@(231, 1.0, 55.82)
    ,@(232, 1.0, 41.02)
    ,@(233, 1.0, 59.82)
    ,@(234, 2.0, 54.89)
    ,@(235, 2.0, 57.23)
    ,@(236, 2.0, 68.71)
    ,@(237, 2.0, 58.3)
    ,@(238, 2.0, 67.82)
    ,@(239, 2.0, 69.51)
    ,@(240, 2.0, 58.68)
    ,@(241, 3.0, 55.07)
    ,@(242, 2.0, 59.96)
    ,@(243, 2.0, 58.23)
    ,@(244, 1.0, 35.97)
    ,@(245, 1.0, 50.41)
    ,@(246, 3.0, 58.65)
    ,@(247, 1.0, 56.57)
    ,@(248, 1.0, 33.09)
    ,@(249, 1.0, 33.09)
    ,@(250, 2.0, 71.45)
    ,@(251, 1.0, 16.57)
    ,@(252, 1.0, 75.96)
    ,@(253, 1.0, 54.17)
    ,@(254, 2.0, 59.36)
    ,@(255, 2.0, 69.54)
    ,@(256, 2.0, 45.45)
    ,@(257, 2.0, 61.16)
    ,@(258, 2.0, 67.81)
    ,@(259, 1.0, 68.97)
    ,@(260, 2.0, 59.9)
    ,@(261, 3.0, 69.06)
    ,@(262, 3.0, 68.08)
    ,@(263, 2.0, 42.34)
    ,@(264, 1.0, 58.25)
    ,@(265, 1.0, 56.41)
    ,@(266, 1.0, 77.94)
    ,@(267, 1.0, 74.05)
    ,@(268, 1.0, 66.08)
    ,@(269, 1.0, 43.54)
    ,@(270, 2.0, 50.44)
    ,@(271, 2.0, 51.75)
    ,@(272, 2.0, 70.96)
    ,@(273, 1.0, 66.29)
    ,@(274, 2.0, 60.42)
    ,@(275, 2.0, 42.34)
    ,@(277, 3.0, 61.95)
    ,@(278, 1.0, 62.51)
    ,@(279, 2.0, 52.66)
    ,@(280, 2.0, 55.31)
    ,@(281, 2.0, 41.51)
    ,@(282, 2.0, 69.9)
    ,@(283, 2.0, 54.32)
    ,@(284, 2.0, 80.92)
    ,@(285, 3.0, 70.53)
    ,@(286, 2.0, 56.85)
    ,@(287, 2.0, 51.48)
    ,@(288, 2.0, 75.05)
    ,@(289, 2.0, 60.82)
    ,@(290, 1.0, 75.87)
    ,@(291, 2.0, 65.83)
    ,@(292, 3.0, 74.39)
    ,@(293, 3.0, 62.47)
    ,@(294, 2.0, 66.28)
    ,@(295, 2.0, 67.42)
    ,@(296, 2.0, 73.12)
)

$ws.Range("E2:F2").Copy()
foreach ($entry in $data) {
    $r = $entry[0]
    $target = $ws.Range("E" + $r + ":F" + $r)
    $target.PasteSpecial(-4122)  # xlPasteFormats: reuse the existing column style
    $ws.Cells.Item($r, 5).Value = $entry[1]
    $ws.Cells.Item($r, 6).Value = $entry[2]
}
$excel.CutCopyMode = $false
